# Facebook registration page test data — add Sheet2 after Sheet1 and
# populate the registration row (name, surname, phone, email, dob, gender).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Deepali"
$ws2.Range("B1").Value = "patil"
$ws2.Range("C1").Value = 9082227391
$ws2.Range("D1").Value = "omsai@123"
$ws2.Range("E1").Value = 2
$ws2.Range("F1").Value = 5
$ws2.Range("G1").Value = 1994
$ws2.Range("H1").Value = "Female"

# Turn the email cell into a mailto hyperlink (matches the style already
# used for the hyperlink cell on Sheet1).
$ws2.Hyperlinks.Add($ws2.Range("D1"), "mailto:omsai@123", "", "", "omsai@123")
$ws1.Range("B1").Copy()
$ws2.Range("D1").PasteSpecial(-4122)
